# ============================================================================
# PlayerPerformance_4387.xlsx edit script
#
# Summary of the change:
#   1. Insert a brand-new "Player Info" worksheet as the FIRST sheet with
#      player bio columns (ID, NAME, BATTING_HAND, BOWL_STYLE).
#   2. On the existing "ODI Batting" sheet: rename the MATCH_CARD_LINK
#      column to MATCH_CODE and replace each URL value with just the
#      numeric match code. Also clear a few stray empty INNING_NUMBER
#      cells (rows where the player did not bat).
#   3. On the existing "ODI Bowling" sheet: same MATCH_CARD_LINK ->
#      MATCH_CODE rename + URL -> code replacement.
#   4. Append a brand-new "ODI Batting Extra" worksheet (after "ODI
#      Bowling") with additional per-match batting detail columns.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Helper: copy the header formatting (bold, centered, thin border) from
# an already-styled header cell onto a destination range, reusing the
# workbook's existing header style instead of fabricating a new one.
# NOTE: re-resolve the template range *each call* (rather than once,
# up-front) -- a Range handle captured before a Worksheets.Add() can end
# up pointing at the wrong sheet once indices shift underneath it.
# ----------------------------------------------------------------------
function Set-HeaderStyle($range) {
    $headerTemplate = $wb.Worksheets.Item("ODI Batting").Range("A1")
    $headerTemplate.Copy()
    $range.PasteSpecial(-4122)
}

# ----------------------------------------------------------------------
# Helper: write a value into a cell while forcing TEXT storage, even
# when the value looks numeric (match codes, percentages, etc. must stay
# text, matching the source scrape). A leading apostrophe is Excel's
# classic "treat as text" quote-prefix; it also works for an empty
# string, producing an empty (but text-typed) cell rather than skipping
# the cell outright.
# ----------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = "'" + $value
}

# ============================================================================
# 1. New "Player Info" sheet (inserted before the active/first sheet)
# ============================================================================
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle($playerInfo.Range("A1:D1"))

Set-TextValue $playerInfo 2 1 "4387"
$playerInfo.Range("B2").Value = "Shreyas Santosh Iyer"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# ============================================================================
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE
# ============================================================================
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2  = "4096"; 3  = "4098"; 4  = "4099"; 5  = "4130"; 6  = "4133"
    7  = "4135"; 8  = "4359"; 9  = "4360"; 10 = "4362"; 11 = "4385"
    12 = "4387"; 13 = "4388"; 14 = "4398"; 15 = "4399"; 16 = "4400"
    17 = "4402"; 18 = "4406"; 19 = "4410"; 20 = "4435"; 21 = "4436"
    22 = "4437"; 23 = "4454"; 24 = "4524"; 25 = "4526"; 26 = "4529"
    27 = "4536"; 28 = "4609"; 29 = "4621"; 30 = "4623"; 31 = "4624"
    32 = "4656"; 33 = "4657"; 34 = "4658"; 35 = "4669"; 36 = "4673"
    37 = "4676"; 38 = "4679"; 39 = "4682"; 40 = "4685"; 41 = "4687"
    42 = "4689"; 43 = "4691"
}
foreach ($row in $battingCodes.Keys) {
    Set-TextValue $batting $row 4 $battingCodes[$row]
}

# Rows where the player did not bat had a stray blank INNING_NUMBER cell;
# clear those so the cell disappears entirely (matches did-not-bat rows).
$emptyInningRows = @(7, 8, 28, 36)
foreach ($row in $emptyInningRows) {
    $batting.Cells.Item($row, 2).ClearContents()
}

# ============================================================================
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE
# ============================================================================
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4098"; 3 = "4387"; 4 = "4526"; 5 = "4529"; 6 = "4691"
}
foreach ($row in $bowlingCodes.Keys) {
    Set-TextValue $bowling $row 2 $bowlingCodes[$row]
}

# ============================================================================
# 4. New "ODI Batting Extra" sheet (appended after "ODI Bowling")
# ============================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle($extra.Range("A1:F1"))

# MATCH_CODE, BATTING_POSITION (numeric or blank), NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4524", "",  "",   "",   "",       "NO"),
    @("4526", "5", "0",  "0",  "3.83%",  "NO"),
    @("4529", "5", "2",  "0",  "9.19%",  "NO"),
    @("4536", "4", "9",  "0",  "30.19%", "YES"),
    @("4609", "3", "",   "",   "",       "NO"),
    @("4621", "3", "5",  "2",  "17.53%", "NO"),
    @("4623", "",  "",   "",   "",       "NO"),
    @("4624", "3", "4",  "1",  "19.56%", "NO"),
    @("4656", "",  "",   "",   "",       "NO"),
    @("4657", "4", "15", "0",  "40.07%", "YES"),
    @("4658", "4", "3",  "2",  "26.67%", "NO"),
    @("4669", "3", "4",  "4",  "26.14%", "NO"),
    @("4673", "",  "",   "",   "",       "NO"),
    @("4676", "3", "8",  "0",  "22.37%", "NO"),
    @("4679", "4", "2",  "0",  "12.90%", "NO"),
    @("4682", "3", "6",  "3",  "30.83%", "NO"),
    @("4685", "",  "",   "",   "",       "NO"),
    @("4687", "",  "",   "",   "",       "NO"),
    @("4689", "4", "5",  "0",  "12.79%", "NO"),
    @("4691", "",  "",   "",   "",       "NO")
)

$r = 2
foreach ($row in $extraRows) {
    Set-TextValue $extra $r 1 $row[0]
    if ($row[1] -ne "") {
        # BATTING_POSITION is a genuine number when present ...
        $extra.Cells.Item($r, 2).Value = [double]$row[1]
    } else {
        # ... but an empty (still text-typed) cell when the player
        # didn't bat in that match.
        Set-TextValue $extra $r 2 $row[1]
    }
    Set-TextValue $extra $r 3 $row[2]
    Set-TextValue $extra $r 4 $row[3]
    Set-TextValue $extra $r 5 $row[4]
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# Land the selection back on the first sheet / A1, matching a freshly
# opened workbook.
$playerInfo.Activate()
$playerInfo.Range("A1").Select()
